# Update "想去人数" (number of people interested) values in the "F" column
# across the relevant worksheets, as published to gh-pages output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet: 展览 (Exhibition)
$ws1.Range("F8").Value = 5136
$ws1.Range("F9").Value = 1443
$ws1.Range("F15").Value = 4131
$ws1.Range("F16").Value = 973
$ws1.Range("F24").Value = 942
$ws1.Range("F25").Value = 280
$ws1.Range("F29").Value = 1055
$ws1.Range("F34").Value = 198
$ws1.Range("F35").Value = 1590
$ws1.Range("F36").Value = 2128
$ws1.Range("F37").Value = 993
$ws1.Range("F38").Value = 31
$ws1.Range("F39").Value = 240
$ws1.Range("F40").Value = 583
$ws1.Range("F41").Value = 245
$ws1.Range("F47").Value = 123

# Sheet: 本地生活 (Local Life)
$ws3.Range("F2").Value = 714

# Sheet: 全部类型 (All Types)
$ws4.Range("F2").Value = 714
$ws4.Range("F9").Value = 5136
$ws4.Range("F10").Value = 1443
$ws4.Range("F16").Value = 4131
$ws4.Range("F17").Value = 973
$ws4.Range("F28").Value = 942
$ws4.Range("F29").Value = 280
$ws4.Range("F33").Value = 1055
$ws4.Range("F35").Value = 1590
$ws4.Range("F36").Value = 2128
$ws4.Range("F38").Value = 993
$ws4.Range("F39").Value = 31
$ws4.Range("F41").Value = 240
$ws4.Range("F42").Value = 583
$ws4.Range("F43").Value = 245
$ws4.Range("F48").Value = 123
